$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Charts")
$ws.Range("A1").Value = "Automatically generated chart(s) coming soon to this tab."
$ws.Activate()
